# Regenerate save_data: column G (header "K" = strikeouts) had been
# populated from the wrong source stat (Strike#). Recalculate and
# rewrite the correct per-start K values for rows 2-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(3, 3, 3, 4, 3, 5, 5, 1, 6, 2, 3, 1, 2, 3, 7, 7, 6, 5, 2, 2, 6, 6, 2, 3, 7, 2, 2, 2, 2)

$firstRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
